# Report_Face_Detection.xlsx update
# - Moves the "* hướng thực hiện" sub-block one row up and one column to the
#   right (B92/C93:C95 -> C91/D92:D94), nesting it deeper under B83's block.
# - Adds two new follow-up blocks (rows 100-101 and 103-104) documenting
#   progress on the "liveness detection" (#4) and "face angle" (#2) items.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "huong thuc hien" sub-block up one row and right one column ---

# Grab the source values before we start overwriting anything.
$vB92 = $ws.Cells.Item(92, 2).Value2
$vC93 = $ws.Cells.Item(93, 3).Value2
$vC94 = $ws.Cells.Item(94, 3).Value2
$vC95 = $ws.Cells.Item(95, 3).Value2

# Write the moved values into their new homes.
$ws.Cells.Item(91, 3).Value2 = $vB92
$ws.Cells.Item(92, 4).Value2 = $vC93
$ws.Cells.Item(93, 4).Value2 = $vC94
$ws.Cells.Item(94, 4).Value2 = $vC95

# Carry over the quotePrefix formatting that C93:C95 had onto D92:D94.
$ws.Cells.Item(93, 3).Copy()
$ws.Cells.Item(92, 4).PasteSpecial(-4122)
$ws.Cells.Item(94, 4).PasteSpecial(-4122)
$ws.Cells.Item(93, 4).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Clear out the old locations now that the content has been relocated
# (use Clear(), not ClearContents(), so the now-unused cells disappear
# entirely rather than leaving a bare formatted cell behind).
$ws.Cells.Item(92, 2).Clear()
$ws.Cells.Item(93, 3).Clear()
$ws.Cells.Item(94, 3).Clear()
$ws.Cells.Item(95, 3).Clear()

# --- Add the new "liveness detection" follow-up block ---
$ws.Cells.Item(100, 2).Value2 = "* Hiện tại phát hiện trong khung hình có sự sống thì cho ảnh vào vẫn phát hiện ra người."
$ws.Cells.Item(101, 3).Value2 = "'- Đang kiểm tra lại source code phần xử lý kiểm tra phát hiện sự sống trong khung hình thì mới tiếp tục nhận diện khuôn mặt."

# --- Add the new "face angle" follow-up block ---
$ws.Cells.Item(103, 2).Value2 = "* Chưa thể phát hiện được những khuôn mặt ở khía độ góc cạnh như nghiêng hay cúi"
$ws.Cells.Item(104, 3).Value2 = "'- Hiện tại sẽ bỏ qua không xử lý, sẽ tập trung xử lý nâng cao độ chính xác và tốc độ xử lý nhận diện với khuôn mặt chính diện"

# --- Update the view so the sheet opens scrolled to the newly edited area ---
$ws.Activate()
try { $excel.ActiveWindow.ScrollRow = 73 } catch {}
$ws.Range("C102").Select()
